# Updated scripts for TSEE Reports
# Adds a new "Global Filters" worksheet (with a "Beta Filters" label and the
# CA / ON / Toronto / Neural Turing Tech location row that used to live on
# Reviews_Filter!C2:F2), and blanks out the corresponding cells on
# Reviews_Filter to "null" now that they live on the new sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new "Global Filters" sheet after the last existing tab ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$gf = $wb.Worksheets.Add($null, $lastSheet)
$gf.Name = "Global Filters"

$gf.Range("A3").Value = "Beta Filters"

$gf.Range("A5").Value = "CA"
$gf.Range("B5").Value = "ON"
$gf.Range("C5").Value = "Toronto"
$gf.Range("D5").Value = "Neural Turing Tech - Primrose, 1131 Steeles Ave. West, M2R 3W8, +14164510870"

$gf.Range("A5:D5").Select() | Out-Null

# --- Update Reviews_Filter row 2 filter values to "null" ---
$reviewsFilter = $wb.Worksheets.Item("Reviews_Filter")
$reviewsFilter.Range("C2").Value = "null"
$reviewsFilter.Range("D2").Value = "null"
$reviewsFilter.Range("E2").Value = "null"
$reviewsFilter.Range("F2").Value = "null"

# Restore Reviews_Filter as the active sheet (unchanged from the original file)
$reviewsFilter.Activate()
